$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (Shelter Assigned, Shelter Latitude, Shelter Longitude)
$changes = @{
    3  = @("San Marcos Elem. Sch.", 120.778807101888, 14.8978852342351)
    5  = @("Mun. Covered Court", 120.765260684344, 14.9142263888355)
    6  = @("San Miguel Evacuation Center", 120.743854227842, 14.9239360863606)
    7  = @("Barangay Hall Bulusan", 120.742556915324, 14.9104237817507)
    8  = @("Calizon Dike", 120.749282551866, 14.9149192538653)
    9  = @("Danga Dike", 120.75016278348, 14.9271290793079)
    10 = @("Calizon Dike", 120.749282551866, 14.9149192538653)
    13 = @("San Jose E.C.", 120.734404217521, 14.8832968907295)
    15 = @("Mun. Covered Court", 120.765260684344, 14.9142263888355)
    16 = @("Meysulao Multipurpose/E.C.", 120.739306695061, 14.9059685215204)
    18 = @("Mun. Covered Court", 120.765260684344, 14.9142263888355)
    19 = @("GMA Kapuso E.C.", 120.796318926682, 14.8957602136898)
    20 = @("San Miguel Evacuation Center", 120.743854227842, 14.9239360863606)
    21 = @("San Miguel Evacuation Center", 120.743854227842, 14.9239360863606)
    22 = @("BMLTC Multi-Purpose Bldg and EC", 120.766768211462, 14.9185376869108)
    23 = @("San Marcos National H.S.", 120.778052198391, 14.8938559117996)
    25 = @("San Jose E.C.", 120.734404217521, 14.8832968907295)
    26 = @("NV9 Multi-Purpose", 120.765656923103, 14.8984661105575)
    27 = @("F. Mendoza Memorial Elem Sch.", 120.767878787289, 14.9176780529243)
    28 = @("GMA Kapuso E.C.", 120.796318926682, 14.8957602136898)
    29 = @("Balungao E.C.", 120.761492937455, 14.9148551401837)
}

foreach ($row in $changes.Keys) {
    $vals = $changes[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
    $ws.Cells.Item($row, 6).Value = $vals[2]
}
